# Updated cryptos list on Sun Jul  9 09:23:54 UTC 2023 with GitHub Actions
#
# Refreshes the scraped Price (column D) / Volume(1h) (column E) figures for
# each coin row, and corrects the Aptos / EnergySwap ranking swap (rows 47-48)
# picked up in this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (never let Excel auto-convert a numeric-
# looking price string, e.g. "1.002", into a Number) -- every Price cell in
# this sheet is stored as text, matching how the scraper originally wrote it.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Price / Volume(1h) refresh ----------------------------------------

$ws.Range("D2").Value = "30.300.65"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("D3").Value = "1.868.05"
$ws.Range("E3").Value = "  +0.11%  "

Set-TextValue $ws.Range("D4") "1.002"
$ws.Range("E4").Value = "  +0.11%  "

Set-TextValue $ws.Range("D5") "234.86"
$ws.Range("E5").Value = "  -1.01%  "

Set-TextValue $ws.Range("D7") "0.4700"
$ws.Range("E7").Value = "  +0.34%  "

Set-TextValue $ws.Range("D8") "0.2868"
$ws.Range("E8").Value = "  -0.02%  "

Set-TextValue $ws.Range("D9") "0.06578"
$ws.Range("E9").Value = "  +0.34%  "

Set-TextValue $ws.Range("D10") "21.58"
$ws.Range("E10").Value = "  -3.23%  "

Set-TextValue $ws.Range("D11") "0.08012"
$ws.Range("E11").Value = "  +1.39%  "

Set-TextValue $ws.Range("D12") "96.80"
$ws.Range("E12").Value = "  -1.18%  "

$ws.Range("D13").Value = "1.871.48"
$ws.Range("E13").Value = "  +0.25%  "

Set-TextValue $ws.Range("D14") "5.111"
$ws.Range("E14").Value = "  -1.36%  "

Set-TextValue $ws.Range("D15") "0.6838"
$ws.Range("E15").Value = "  +0.23%  "

Set-TextValue $ws.Range("D16") "269.45"
$ws.Range("E16").Value = "  -3.43%  "

$ws.Range("D17").Value = "30.327.90"
$ws.Range("E17").Value = "  +0.15%  "

Set-TextValue $ws.Range("D18") "13.99"
$ws.Range("E18").Value = "  +2.40%  "

Set-TextValue $ws.Range("D19") "0.000007629"
$ws.Range("E19").Value = "  +3.81%  "

Set-TextValue $ws.Range("D20") "1.001"
$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").Value = "2.118.00"
$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("E22").Value = "  +0.07%  "

Set-TextValue $ws.Range("D23") "5.271"
$ws.Range("E23").Value = "  -2.43%  "

Set-TextValue $ws.Range("D24") "6.207"
$ws.Range("E24").Value = "  -0.01%  "

Set-TextValue $ws.Range("D25") "9.402"
$ws.Range("E25").Value = "  +1.12%  "

Set-TextValue $ws.Range("D26") "168.04"
$ws.Range("E26").Value = "  -0.39%  "

Set-TextValue $ws.Range("D27") "18.90"
$ws.Range("E27").Value = "  -1.32%  "

Set-TextValue $ws.Range("D28") "1.948"
$ws.Range("E28").Value = "  -0.06%  "

Set-TextValue $ws.Range("D29") "1.370"
$ws.Range("E29").Value = "  -0.87%  "

Set-TextValue $ws.Range("D30") "0.09915"
$ws.Range("E30").Value = "  +0.65%  "

Set-TextValue $ws.Range("D31") "4.368"
$ws.Range("E31").Value = "  -0.72%  "

$ws.Range("E32").Value = "  -1.44%  "

Set-TextValue $ws.Range("D33") "4.064"
$ws.Range("E33").Value = "  -0.25%  "

Set-TextValue $ws.Range("D34") "0.04717"
$ws.Range("E34").Value = "  -0.74%  "

Set-TextValue $ws.Range("D35") "1.136"
$ws.Range("E35").Value = "  -0.97%  "

Set-TextValue $ws.Range("D36") "0.6999"

Set-TextValue $ws.Range("D37") "2.712"
$ws.Range("E37").Value = "  +0.17%  "

Set-TextValue $ws.Range("D38") "0.01869"
$ws.Range("E38").Value = "  -0.59%  "

Set-TextValue $ws.Range("D39") "2.648"
$ws.Range("E39").Value = "  +1.03%  "

Set-TextValue $ws.Range("D40") "6.274"

Set-TextValue $ws.Range("D41") "71.80"
$ws.Range("E41").Value = "  -6.71%  "

Set-TextValue $ws.Range("D42") "1.957"
$ws.Range("E42").Value = "  -0.41%  "

Set-TextValue $ws.Range("D43") "0.8421"
$ws.Range("E43").Value = "  -1.27%  "

Set-TextValue $ws.Range("D44") "0.4164"
$ws.Range("E44").Value = "  -0.70%  "

Set-TextValue $ws.Range("D45") "0.9996"
$ws.Range("E45").Value = "  +0.01%  "

Set-TextValue $ws.Range("D46") "102.74"
$ws.Range("E46").Value = "  -0.71%  "

# --- Rows 47-48: EnergySwap and Aptos swap ranking places --------------
# Row 47 was EnergySwap; this run it is Aptos (with refreshed figures).
Set-TextValue $ws.Range("D47") "7.044"
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("E47").Value = "  -2.56%  "

# Row 48 was Aptos; this run it is EnergySwap (with refreshed figures).
Set-TextValue $ws.Range("D48") "9.110"
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E48").Value = "  -2.14%  "

Set-TextValue $ws.Range("D49") "905.97"
$ws.Range("E49").Value = "  -5.86%  "

Set-TextValue $ws.Range("D50") "34.45"
$ws.Range("E50").Value = "  +0.37%  "

Set-TextValue $ws.Range("D51") "0.05707"
$ws.Range("E51").Value = "  +1.09%  "
